$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3036.4119
$ws.Range("I137").Value = 2953.2964
$ws.Range("J137").Value = 3357
$ws.Range("K137").Value = 8859.889200000001
$ws.Range("L137").Value = 10071
$ws.Range("M137").Value = -6309.889200000001
$ws.Range("N137").Value = -15171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2930.3125
$ws.Range("I2").Value = 2079.9167
$ws.Range("K2").Value = 2079.9167
$ws.Range("M2").Value = -1966.9167
$ws.Range("H45").Value = 1406.8
$ws.Range("I45").Value = 1273.3334
$ws.Range("J45").Value = 1607
$ws.Range("K45").Value = 1273.3334
$ws.Range("L45").Value = 1607
$ws.Range("M45").Value = -896.3334
$ws.Range("N45").Value = -2361
$ws.Range("H61").Value = 811497.1
$ws.Range("I61").Value = 849008.6
$ws.Range("J61").Value = 773985.7
$ws.Range("K61").Value = 849008.6
$ws.Range("L61").Value = 773985.7
$ws.Range("M61").Value = -848796.6
$ws.Range("N61").Value = -774409.7
$ws.Range("H74").Value = 263919.84
$ws.Range("I74").Value = 304167.28
$ws.Range("J74").Value = 116346
$ws.Range("K74").Value = 304167.28
$ws.Range("L74").Value = 116346
$ws.Range("M74").Value = -303293.28
$ws.Range("N74").Value = -118094
$ws.Range("H77").Value = 263919.84
$ws.Range("I77").Value = 304167.28
$ws.Range("J77").Value = 116346
$ws.Range("K77").Value = 1520836.4
$ws.Range("L77").Value = 581730
$ws.Range("M77").Value = -1516468.4
$ws.Range("N77").Value = -590466
$ws.Range("H88").Value = 2230.7104
$ws.Range("I88").Value = 2374
$ws.Range("J88").Value = 2071.5
$ws.Range("K88").Value = 2374
$ws.Range("L88").Value = 2071.5
$ws.Range("M88").Value = -1968
$ws.Range("N88").Value = -2883.5
$ws.Range("H91").Value = 2230.7104
$ws.Range("I91").Value = 2374
$ws.Range("J91").Value = 2071.5
$ws.Range("K91").Value = 2374
$ws.Range("L91").Value = 2071.5
$ws.Range("M91").Value = -970
$ws.Range("N91").Value = -4879.5
$ws.Range("H97").Value = 809.06665
$ws.Range("I97").Value = 846.5833
$ws.Range("J97").Value = 659
$ws.Range("K97").Value = 846.5833
$ws.Range("L97").Value = 659
$ws.Range("M97").Value = -350.5833
$ws.Range("N97").Value = -1651
$ws.Range("H116").Value = 2930.3125
$ws.Range("I116").Value = 2079.9167
$ws.Range("K116").Value = 2079.9167
$ws.Range("M116").Value = 214.0832999999998
$ws.Range("H122").Value = 3390.7908
$ws.Range("I122").Value = 3271.353
$ws.Range("K122").Value = 9814.059000000001
$ws.Range("M122").Value = -7364.059000000001
$ws.Range("H136").Value = 811497.1
$ws.Range("I136").Value = 849008.6
$ws.Range("J136").Value = 773985.7
$ws.Range("K136").Value = 2547025.8
$ws.Range("L136").Value = 2321957.1
$ws.Range("M136").Value = -2544475.8
$ws.Range("N136").Value = -2327057.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2930.3125
$ws.Range("I3").Value = 2079.9167
$ws.Range("K3").Value = 2079.9167
$ws.Range("M3").Value = -1965.9167
$ws.Range("H107").Value = 2986.75
$ws.Range("I107").Value = 2567
$ws.Range("J107").Value = 3313.2222
$ws.Range("K107").Value = 2567
$ws.Range("L107").Value = 3313.2222
$ws.Range("M107").Value = -647
$ws.Range("N107").Value = -7153.2222

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 308.34784
$ws.Range("I22").Value = 315.53845
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 315.53845
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = 34.46154999999999
$ws.Range("N22").Value = -999
$ws.Range("H58").Value = 4002
$ws.Range("I58").Value = 4888.56
$ws.Range("J58").Value = 2418.8572
$ws.Range("K58").Value = 4888.56
$ws.Range("L58").Value = 2418.8572
$ws.Range("M58").Value = -4685.56
$ws.Range("N58").Value = -2824.8572
$ws.Range("H122").Value = 1619.5625
$ws.Range("I122").Value = 914.25
$ws.Range("J122").Value = 2324.875
$ws.Range("K122").Value = 2742.75
$ws.Range("L122").Value = 6974.625
$ws.Range("M122").Value = -292.75
$ws.Range("N122").Value = -11874.625
$ws.Range("H132").Value = 3017.84
$ws.Range("I132").Value = 1840.5
$ws.Range("J132").Value = 5110.8887
$ws.Range("K132").Value = 5521.5
$ws.Range("L132").Value = 15332.6661
$ws.Range("M132").Value = -2991.5
$ws.Range("N132").Value = -20392.6661
$ws.Range("H136").Value = 4002
$ws.Range("I136").Value = 4888.56
$ws.Range("J136").Value = 2418.8572
$ws.Range("K136").Value = 14665.68
$ws.Range("L136").Value = 7256.571599999999
$ws.Range("M136").Value = -12115.68
$ws.Range("N136").Value = -12356.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 390.41025
$ws.Range("I5").Value = 336.83334
$ws.Range("K5").Value = 1010.50002
$ws.Range("M5").Value = -898.5000200000001
$ws.Range("H51").Value = 300
$ws.Range("I51").Value = 300
$ws.Range("K51").Value = 900
$ws.Range("M51").Value = -440
$ws.Range("H58").Value = 2888
$ws.Range("J58").Value = 3097.7778
$ws.Range("L58").Value = 9293.3334
$ws.Range("N58").Value = -9549.3334
$ws.Range("H64").Value = 2082.6667
$ws.Range("I64").Value = 1188.8
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 3566.4
$ws.Range("L64").Value = 9600
$ws.Range("M64").Value = -3296.4
$ws.Range("N64").Value = -10140
$ws.Range("H67").Value = 2082.6667
$ws.Range("I67").Value = 1188.8
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 3566.4
$ws.Range("L67").Value = 9600
$ws.Range("M67").Value = -2630.4
$ws.Range("N67").Value = -11472
$ws.Range("H122").Value = 810.93335
$ws.Range("I122").Value = 440
$ws.Range("J122").Value = 1831
$ws.Range("K122").Value = 3960
$ws.Range("L122").Value = 16479
$ws.Range("M122").Value = -1510
$ws.Range("N122").Value = -21379
$ws.Range("H135").Value = 390.41025
$ws.Range("I135").Value = 336.83334
$ws.Range("K135").Value = 3031.50006
$ws.Range("M135").Value = -496.5000600000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5926.2104
$ws.Range("I80").Value = 7530.615
$ws.Range("J80").Value = 2450
$ws.Range("K80").Value = 7530.615
$ws.Range("L80").Value = 2450
$ws.Range("M80").Value = -6532.615
$ws.Range("N80").Value = -4446
$ws.Range("H83").Value = 5926.2104
$ws.Range("I83").Value = 7530.615
$ws.Range("J83").Value = 2450
$ws.Range("K83").Value = 37653.075
$ws.Range("L83").Value = 12250
$ws.Range("M83").Value = -32661.075
$ws.Range("N83").Value = -22234
$ws.Range("H102").Value = 6488.615
$ws.Range("I102").Value = 3157.2
$ws.Range("J102").Value = 17593.334
$ws.Range("K102").Value = 3157.2
$ws.Range("L102").Value = 17593.334
$ws.Range("M102").Value = -1535.2
$ws.Range("N102").Value = -20837.334
$ws.Range("H107").Value = 21743724
$ws.Range("I107").Value = 8470.583000000001
$ws.Range("J107").Value = 45454910
$ws.Range("K107").Value = 8470.583000000001
$ws.Range("L107").Value = 45454910
$ws.Range("M107").Value = -6550.583000000001
$ws.Range("N107").Value = -45458750
$ws.Range("H113").Value = 5624.522
$ws.Range("I113").Value = 7135.8125
$ws.Range("K113").Value = 7135.8125
$ws.Range("M113").Value = -4965.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1443.5358
$ws.Range("I16").Value = 1524.7142
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1524.7142
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -1354.7142
$ws.Range("N16").Value = -1540
$ws.Range("H46").Value = 993.2727
$ws.Range("I46").Value = 1074.2858
$ws.Range("J46").Value = 851.5
$ws.Range("K46").Value = 1074.2858
$ws.Range("L46").Value = 851.5
$ws.Range("M46").Value = -886.2858000000001
$ws.Range("N46").Value = -1227.5
$ws.Range("H132").Value = 12206.305
$ws.Range("I132").Value = 5050.375
$ws.Range("J132").Value = 16022.8
$ws.Range("K132").Value = 15151.125
$ws.Range("L132").Value = 48068.39999999999
$ws.Range("M132").Value = -12621.125
$ws.Range("N132").Value = -53128.39999999999
